# Repository-name change in the resume's GitHub Pages hyperlink:
#   "gordon.github.io" -> "blursotongkia.github.io"
#
# The original paragraph is a single <w:hyperlink> containing one run
# whose text is "gordon.github.io". The target state keeps that same
# run for the new "blursotongkia" text, and adds a *second* run (same
# character formatting, i.e. the Hyperlink style/font/color/size) that
# holds ".github.io" - so the run gets split in two instead of just
# having its text swapped in place.

$d = $word.ActiveDocument

# 1) Locate the hyperlink and rewrite its visible text in one shot so
#    the whole "gordon.github.io" -> "blursotongkia.github.io" swap
#    happens as a simple find/replace (this keeps it inside the
#    hyperlink and preserves the run's formatting).
$hl = $d.Hyperlinks.Item(1)
$hlRange = $hl.Range
$found = $hlRange.Find.Execute("gordon.github.io", $true, $false, $false, $false, $false, $true, 1, $false, "blursotongkia.github.io", 2)
if (-not $found) {
    throw "Could not find 'gordon.github.io' in the hyperlink text."
}

# 2) Re-find just the ".github.io" tail of that same text and nudge its
#    character formatting (toggle Bold on, then back off) - this forces
#    the engine to split the run at that boundary while leaving the
#    net formatting identical to its neighbour, giving two runs with
#    matching <w:rPr> (rStyle Hyperlink, fonts, color, sz/szCs) exactly
#    like the rest of the hyperlink.
$tail = $d.Content
$tailFound = $tail.Find.Execute(".github.io", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $tailFound) {
    throw "Could not find '.github.io' tail after replacement."
}
$tail.Font.Bold = 1
$tail.Font.Bold = 0
